# Add two new columns, I ("I0") and J ("IF"), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold, bordered, centered) by
# copying the format from the adjacent header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-47 ------------------------------------------------------
$iValues = @{
    2 = 5;  3 = 7;  4 = 4;  5 = 9;  6 = 6;  7 = 6;  8 = 4;  9 = 5;  10 = 6;
    11 = 9; 12 = 8; 13 = 7; 14 = 9; 15 = 7; 16 = 7; 17 = 7; 18 = 9; 19 = 8;
    20 = 9; 21 = 9; 22 = 8; 23 = 8; 24 = 6; 25 = 8; 26 = 7; 27 = 6; 28 = 4;
    29 = 10; 30 = 9; 31 = 8; 32 = 8; 33 = 8; 34 = 11; 35 = 9; 36 = 7; 37 = 7;
    38 = 8; 39 = 7; 40 = 9; 41 = 8; 42 = 7; 43 = 5; 44 = 7; 45 = 6; 46 = 7;
    47 = 6
}
$jValues = @{
    2 = 5;  3 = 7;  4 = 4;  5 = 9;  6 = 6;  7 = 6;  8 = 4;  9 = 5;  10 = 6;
    11 = 9; 12 = 8; 13 = 7; 14 = 9; 15 = 7; 16 = 7; 17 = 7; 18 = 9; 19 = 8;
    20 = 9; 21 = 9; 22 = 8; 23 = 8; 24 = 7; 25 = 8; 26 = 7; 27 = 6; 28 = 4;
    29 = 10; 30 = 9; 31 = 8; 32 = 8; 33 = 8; 34 = 11; 35 = 9; 36 = 7; 37 = 7;
    38 = 8; 39 = 7; 40 = 9; 41 = 8; 42 = 7; 43 = 5; 44 = 7; 45 = 6; 46 = 7;
    47 = 6
}

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
